$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, preventing Excel from auto-coercing
# numeric-looking strings (e.g. "568.65") into real numbers, and leaving the
# cell style untouched (reset to Normal/General afterwards) to match the
# original "General" formatted, string-typed cells.
function Set-CellText($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# --- Simple Price (D) / Volume(1h) (E) updates ---
Set-CellText "D2" "69.536.05"
Set-CellText "E2" "  -1.24%  "
Set-CellText "D3" "2.490.56"
Set-CellText "E3" "  -1.91%  "
Set-CellText "E4" "  +0.13%  "
Set-CellText "D5" "568.65"
Set-CellText "E5" "  -2.34%  "
Set-CellText "D6" "165.36"
Set-CellText "E6" "  -2.57%  "
Set-CellText "E7" "  +0.07%  "
Set-CellText "D8" "0.511"
Set-CellText "E8" "  -1.69%  "
Set-CellText "D9" "2.488.05"
Set-CellText "E9" "  -2.02%  "
Set-CellText "D10" "0.157"
Set-CellText "E10" "  -3.93%  "
Set-CellText "E11" "  -0.51%  "
Set-CellText "D12" "0.353"
Set-CellText "E12" "  +0.03%  "
Set-CellText "D13" "4.92"
Set-CellText "E13" "  +0.41%  "
Set-CellText "D14" "2.951.22"
Set-CellText "E14" "  -1.21%  "
Set-CellText "D15" "69.450.16"
Set-CellText "E15" "  -1.12%  "
Set-CellText "D16" "0.0000174"
Set-CellText "E16" "  -1.71%  "
Set-CellText "D17" "24.48"
Set-CellText "E17" "  -3.49%  "
Set-CellText "D18" "2.499.98"
Set-CellText "E18" "  -1.76%  "
Set-CellText "D19" "11.14"
Set-CellText "E19" "  -2.67%  "
Set-CellText "D20" "7.37"
Set-CellText "E20" "  -7.05%  "
Set-CellText "D21" "345.97"
Set-CellText "E21" "  -2.31%  "
Set-CellText "D22" "3.87"
Set-CellText "E22" "  -2.44%  "
Set-CellText "D23" "1.91"
Set-CellText "E23" "  -3.81%  "
Set-CellText "E24" "  -0.07%  "
Set-CellText "D25" "70.36"
Set-CellText "E25" "  +0.42%  "
Set-CellText "D26" "3.87"
Set-CellText "E26" "  -4.08%  "
Set-CellText "D27" "2.623.29"
Set-CellText "E27" "  -2.03%  "
Set-CellText "D28" "8.59"
Set-CellText "E28" "  -5.72%  "
Set-CellText "E29" "  +0.54%  "
Set-CellText "D32" "451.75"
Set-CellText "E32" "  -3.81%  "
Set-CellText "D33" "1.19"
Set-CellText "E33" "  -6.55%  "
Set-CellText "D38" "19.04"
Set-CellText "E38" "  +0.09%  "
Set-CellText "D39" "18.24"
Set-CellText "E39" "  -1.90%  "
Set-CellText "D41" "0.315"
Set-CellText "E41" "  -2.51%  "
Set-CellText "D42" "4.63"
Set-CellText "E42" "  -4.09%  "
Set-CellText "D43" "1.58"
Set-CellText "E43" "  -2.09%  "
Set-CellText "D44" "38.02"
Set-CellText "E44" "  -0.72%  "
Set-CellText "D45" "2.15"
Set-CellText "E45" "  -8.30%  "
Set-CellText "E46" "  -8.90%  "
Set-CellText "D47" "139.82"
Set-CellText "E47" "  -3.33%  "
Set-CellText "D48" "3.45"
Set-CellText "E48" "  -2.35%  "
Set-CellText "D49" "0.512"
Set-CellText "E49" "  -4.33%  "
Set-CellText "E50" "  -1.07%  "
Set-CellText "D51" "0.572"
Set-CellText "E51" "  -2.33%  "

# --- Row reorderings (rank swaps) with updated B/C/D/E values ---
Set-CellText "B30" "InternetComputer(DFINITY)"
Set-CellText "C30" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-CellText "D30" "7.81"
Set-CellText "E30" "  -1.30%  "
Set-CellText "B31" "PEPE"
Set-CellText "C31" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-CellText "D31" "0.0₃0875"
Set-CellText "E31" "  -4.73%  "
Set-CellText "B34" "FirstDigitalUSD"
Set-CellText "C34" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-CellText "D34" "1.00"
Set-CellText "E34" "  +0.16%  "
Set-CellText "B35" "PancakeSwap"
Set-CellText "C35" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-CellText "D35" "1.71"
Set-CellText "E35" "  -3.07%  "
Set-CellText "B36" "Monero"
Set-CellText "C36" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-CellText "D36" "155.26"
Set-CellText "E36" "  +0.19%  "
Set-CellText "B37" "Kaspa"
Set-CellText "C37" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-CellText "D37" "0.114"
Set-CellText "E37" "  -4.86%  "
